# Update dispatch sheet: insert a new "Field 1 - 50" dispatch row, and set
# the dispatcher ("折篁忆文"/zhuzi) for the existing "Field 51 - 55" row.
#
# The underlying sheet keeps a blank "separator" row above every block of
# assigned rows (row 11 was the blank separator before the "Field 61-65"
# block). We shift that block down by one row (rows 12-14 -> 13-15), turn
# the vacated row 12 into the new "Field 51 - 55" / zhuzi row, and create a
# brand new blank separator row at row 11 for "Field 1 - 50".
# Finally, a new blank row is appended at the bottom (row 23) to keep the
# sheet's trailing blank rows intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append a new blank row 23 (copy of the standard blank template row) ---
$ws.Range("A22:G22").Copy()
$ws.Range("A23:G23").PasteSpecial(-4122)
$ws.Rows(23).RowHeight = 20.35
$ws.Range("A23").Value2 = $ws.Range("A22").Value2
$ws.Range("B23").Value2 = $ws.Range("B22").Value2
$ws.Range("C23").Value2 = $ws.Range("C22").Value2

# --- 2. Shift the "Field 61-65 / Field 66-70 / Field 87-88" block down by one row ---
# row14 -> row15
$ws.Range("A14:C14").Copy()
$ws.Range("A15:C15").PasteSpecial(-4122)
$ws.Range("A15").Value2 = $ws.Range("A14").Value2
$ws.Range("B15").Value2 = $ws.Range("B14").Value2
$ws.Range("C15").Value2 = $ws.Range("C14").Value2

# row13 -> row14
$ws.Range("A13:C13").Copy()
$ws.Range("A14:C14").PasteSpecial(-4122)
$ws.Range("A14").Value2 = $ws.Range("A13").Value2
$ws.Range("B14").Value2 = $ws.Range("B13").Value2
$ws.Range("C14").Value2 = $ws.Range("C13").Value2

# row12 -> row13
$ws.Range("A12:C12").Copy()
$ws.Range("A13:C13").PasteSpecial(-4122)
$ws.Range("A13").Value2 = $ws.Range("A12").Value2
$ws.Range("B13").Value2 = $ws.Range("B12").Value2
$ws.Range("C13").Value2 = $ws.Range("C12").Value2

# --- 3. Turn row 12 into the "Field 51 - 55" row, dispatched to zhuzi (折篁忆文) ---
# B12 needs the "name" cell style (same as B3 / B9); C12 needs the blank
# "dropdown" style (same as the separator row's C column, e.g. C7).
$ws.Range("B9:B9").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("C7:C7").Copy()
$ws.Range("C12").PasteSpecial(-4122)

$ws.Range("A12").Value2 = "Field 51 - 55"
$ws.Range("B12").Value2 = "折篁忆文"
$ws.Range("C12").ClearContents()

# --- 4. Turn row 11 into the new blank separator row for "Field 1 - 50" ---
# A11 needs the normal "field name" style (same as the other A-column field
# cells), B11 needs the style used for dispatcher names on a separator row
# (same as B7), C11 needs the plain blank style (same as C10).
$ws.Range("A13:A13").Copy()
$ws.Range("A11").PasteSpecial(-4122)
$ws.Range("B7:B7").Copy()
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("C10:C10").Copy()
$ws.Range("C11").PasteSpecial(-4122)

$ws.Range("A11").Value2 = "Field 1 - 50"
$ws.Range("B11").Value2 = "小智分配"
$ws.Range("C11").ClearContents()
